# Update reaction sensitivity values on both sheets (NBR, BAR)
$wb = $excel.ActiveWorkbook

$wsNBR = $wb.Worksheets.Item("NBR")
$wsBAR = $wb.Worksheets.Item("BAR")

# New "Reaction_number" (column C) values for rows 2-20 on the NBR sheet
$nbrValues = @(564, 561, 549, 546, 539, 527, 517, 518, 510, 505, 503, 499, 497, 495, 489, 484, 482, 473, 466)

# New "Reaction_number" (column C) values for rows 2-20 on the BAR sheet
$barValues = @(652, 648, 660, 654, 652, 653, 651, 650, 651, 650, 640, 639, 638, 636, 634, 632, 630, 621, 625)

for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
}

for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
